$d = $word.ActiveDocument

# --- Step 1 -----------------------------------------------------------
# The last paragraph in the document (the "y - A discount factor..." bullet)
# gains a <w:rStyle w:val="mi"/> on its paragraph-mark run properties
# (w:pPr/w:rPr), and loses the _GoBack bookmark (it is relocated onto the
# new paragraph added in step 2). Replacing the whole paragraph via
# InsertXML is the most reliable way to rewrite its w:pPr/w:rPr in one shot.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertXML('<w:p w:rsidR="00BA38E7" w:rsidRPr="005C5021" w:rsidRDefault="005C5021" w:rsidP="00BA38E7"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="center" w:pos="3675"/></w:tabs><w:spacing w:after="197" w:line="259" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rStyle w:val="mi"/><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:sz w:val="20"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="005C5021"><w:rPr><w:rStyle w:val="mi"/><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>γ</w:t></w:r><w:r><w:rPr><w:rStyle w:val="mi"/><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:rStyle w:val="mi"/><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:iCs/><w:szCs w:val="29"/></w:rPr><w:t>A discount factor. The closer to one it is, the less it encourages the agent to reach rewards as quickly as possible.</w:t></w:r></w:p>')

# --- Step 2 -----------------------------------------------------------
# Append a brand-new paragraph with the "Policy" prose and the optimal
# value-function formula. Character-style references (rStyle="mi"/"mo")
# on individual math runs do not survive a raw InsertXML insert in this
# host, so the paragraph is inserted with plain run formatting first and
# the "mi"/"mo" character styles are then (re)applied run-by-run below,
# exactly as typing + Format>Style would do it.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p><w:pPr><w:tabs><w:tab w:val="center" w:pos="3675"/></w:tabs><w:spacing w:after="197" w:line="259" w:lineRule="auto"/><w:ind w:left="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Policy </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">of an agent is the strategy or set of rules it utilizes in order to decide what actions it should take at each state in the environment. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">The optimal policy gives us the value function </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>V</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>∗</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>)=max</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>∈</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Caligraphic" w:hAnsi="MathJax_Caligraphic"/><w:b/></w:rPr><w:t>A</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Size1" w:hAnsi="MathJax_Size1"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>[</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>R</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>)+</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>γ</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Size1" w:hAnsi="MathJax_Size1"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>∑</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="15"/><w:szCs w:val="15"/></w:rPr><w:t>′</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>∈</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Caligraphic" w:hAnsi="MathJax_Caligraphic"/><w:b/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>P</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>)(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>′</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>V</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>∗</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Math" w:hAnsi="MathJax_Math"/><w:b/><w:i/><w:iCs/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>′</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Main" w:hAnsi="MathJax_Main"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="MathJax_Size1" w:hAnsi="MathJax_Size1"/><w:b/><w:sz w:val="29"/><w:szCs w:val="29"/></w:rPr><w:t>]</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>')

# --- Step 3 -----------------------------------------------------------
# Re-apply the "mi"/"mo" MathJax character styles onto the individual
# formula runs of the paragraph just inserted.
$newPara = $d.Paragraphs.Last
$base = $newPara.Range.Start
$d.Range($base + 195, $base + 196).Style = "mi"
$d.Range($base + 196, $base + 197).Style = "mo"
$d.Range($base + 197, $base + 198).Style = "mo"
$d.Range($base + 198, $base + 199).Style = "mi"
$d.Range($base + 199, $base + 204).Style = "mo"
$d.Range($base + 204, $base + 205).Style = "mi"
$d.Range($base + 205, $base + 206).Style = "mo"
$d.Range($base + 206, $base + 207).Style = "mi"
$d.Range($base + 207, $base + 208).Style = "mo"
$d.Range($base + 208, $base + 209).Style = "mi"
$d.Range($base + 209, $base + 210).Style = "mo"
$d.Range($base + 210, $base + 211).Style = "mi"
$d.Range($base + 211, $base + 212).Style = "mo"
$d.Range($base + 212, $base + 213).Style = "mi"
$d.Range($base + 213, $base + 215).Style = "mo"
$d.Range($base + 215, $base + 216).Style = "mi"
$d.Range($base + 216, $base + 217).Style = "mo"
$d.Range($base + 217, $base + 218).Style = "mi"
$d.Range($base + 218, $base + 219).Style = "mo"
$d.Range($base + 219, $base + 220).Style = "mo"
$d.Range($base + 220, $base + 221).Style = "mi"
$d.Range($base + 221, $base + 222).Style = "mi"
$d.Range($base + 222, $base + 223).Style = "mo"
$d.Range($base + 223, $base + 224).Style = "mi"
$d.Range($base + 224, $base + 225).Style = "mo"
$d.Range($base + 225, $base + 226).Style = "mi"
$d.Range($base + 226, $base + 228).Style = "mo"
$d.Range($base + 228, $base + 229).Style = "mi"
$d.Range($base + 229, $base + 230).Style = "mo"
$d.Range($base + 230, $base + 231).Style = "mo"
$d.Range($base + 231, $base + 232).Style = "mi"
$d.Range($base + 232, $base + 233).Style = "mo"
$d.Range($base + 233, $base + 234).Style = "mo"
$d.Range($base + 234, $base + 235).Style = "mi"
$d.Range($base + 235, $base + 236).Style = "mo"
$d.Range($base + 236, $base + 237).Style = "mo"
$d.Range($base + 237, $base + 238).Style = "mo"

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
